# Append newly curated gene-list entries (Berlin / Desden panels, 2025-06-19)
# and highlight duplicate gene symbols in column B.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A14").Value = 13
$ws.Range("B14").Value = "ATP7B"
$ws.Range("C14").Value = "Wilson disease (OMIM #277900; frequent diagnostics in Berlin)"
$ws.Range("D14").Value = "yes"
$ws.Range("E14").Value = "2025-06-19"

$ws.Range("A15").Value = 14
$ws.Range("B15").Value = "PAH"
$ws.Range("C15").Value = "Phenylketonuria (OMIM #261600; frequent diagnostics in Berlin)"
$ws.Range("D15").Value = "yes"
$ws.Range("E15").Value = "2025-06-19"

$ws.Range("A16").Value = 15
$ws.Range("B16").Value = "G6PC"
$ws.Range("C16").Value = "Glycogen storage disease Ia (OMIM #232200; frequent diagnostics in Berlin)"
$ws.Range("D16").Value = "yes"
$ws.Range("E16").Value = "2025-06-19"

$ws.Range("A17").Value = 16
$ws.Range("B17").Value = "HEXA"
$ws.Range("C17").Value = "Tay-Sachs disease (OMIM #272800; frequent diagnostics in Berlin)"
$ws.Range("D17").Value = "yes"
$ws.Range("E17").Value = "2025-06-19"

$ws.Range("A18").Value = 17
$ws.Range("B18").Value = "ACVRL1"
$ws.Range("C18").Value = "Telangiectasia, hereditary hemorrhagic, type 2 (OMIM #600376; frequent diagnostics in Berlin)"
$ws.Range("D18").Value = "yes"
$ws.Range("E18").Value = "2025-06-19"

$ws.Range("A19").Value = 18
$ws.Range("B19").Value = "ENG"
$ws.Range("C19").Value = "Telangiectasia, hereditary hemorrhagic, type 1 (OMIM #187300; frequent diagnostics in Berlin)"
$ws.Range("D19").Value = "yes"
$ws.Range("E19").Value = "2025-06-19"

$ws.Range("A20").Value = 19
$ws.Range("B20").Value = "NOTCH3"
$ws.Range("C20").Value = "Cerebral arteriopathy with subcortical infarcts and leukoencephalopathy 1 (OMIM #125310; frequent diagnostics in Berlin)"
$ws.Range("D20").Value = "yes"
$ws.Range("E20").Value = "2025-06-19"

$ws.Range("A21").Value = 20
$ws.Range("B21").Value = "TREX1"
$ws.Range("C21").Value = "Aicardi-Goutieres syndrome 1, dominant and recessive (OMIM #225750; frequent diagnostics in Berlin)"
$ws.Range("D21").Value = "yes"
$ws.Range("E21").Value = "2025-06-19"

$ws.Range("A22").Value = 21
$ws.Range("B22").Value = "KCNE1"
$ws.Range("C22").Value = "Long QT syndrome 5 (OMIM #613695; frequent diagnostics in Berlin)"
$ws.Range("D22").Value = "yes"
$ws.Range("E22").Value = "2025-06-19"

$ws.Range("A23").Value = 22
$ws.Range("B23").Value = "KCNQ1"
$ws.Range("C23").Value = "Long QT syndrome 1 (OMIM #192500; frequent diagnostics in Berlin)"
$ws.Range("D23").Value = "yes"
$ws.Range("E23").Value = "2025-06-19"

$ws.Range("A24").Value = 23
$ws.Range("B24").Value = "LEMD3"
$ws.Range("C24").Value = "Osteopoikilosis with or without melorheostosis (OMIM #166700; frequent diagnostics in Berlin)"
$ws.Range("D24").Value = "yes"
$ws.Range("E24").Value = "2025-06-19"

$ws.Range("A25").Value = 24
$ws.Range("B25").Value = "UBA1"
$ws.Range("C25").Value = "VEXAS syndrome, somatic (OMIM #`t301054; frequent diagnostics in Berlin)"
$ws.Range("D25").Value = "yes"
$ws.Range("E25").Value = "2025-06-19"

$ws.Range("A26").Value = 25
$ws.Range("B26").Value = "CASR"
$ws.Range("C26").Value = "Hypocalciuric hypercalcemia, type I (OMIM #`t145980; frequent diagnostics in Berlin)"
$ws.Range("D26").Value = "yes"
$ws.Range("E26").Value = "2025-06-19"

$ws.Range("A27").Value = 26
$ws.Range("B27").Value = "ALPL"
$ws.Range("C27").Value = "Hypophosphatasia, adult (OMIM #146300; frequent diagnostics in Berlin)"
$ws.Range("D27").Value = "yes"
$ws.Range("E27").Value = "2025-06-19"

$ws.Range("A28").Value = 27
$ws.Range("B28").Value = "TSHR"
$ws.Range("C28").Value = "Hyperthyroidism, nonautoimmune (OMIM #609152; frequent diagnostics in Berlin)"
$ws.Range("D28").Value = "yes"
$ws.Range("E28").Value = "2025-06-19"

$ws.Range("A29").Value = 28
$ws.Range("B29").Value = "MEFV"
$ws.Range("C29").Value = "Familial Mediterranean fever, AR (OMIM #249100; frequent diagnostics in Berlin)"
$ws.Range("D29").Value = "yes"
$ws.Range("E29").Value = "2025-06-19"

$ws.Range("A30").Value = 29
$ws.Range("B30").Value = "C9ORF72"
$ws.Range("C30").Value = "Frontotemporal dementia and/or amyotrophic lateral sclerosis 1 (OMIM #105550; frequent diagnostics in Berlin)"
$ws.Range("D30").Value = "yes"
$ws.Range("E30").Value = "2025-06-19"

$ws.Range("A31").Value = 30
$ws.Range("B31").Value = "FUS"
$ws.Range("C31").Value = "Amyotrophic lateral sclerosis 6, with or without frontotemporal dementia (OMIM #608030; frequent diagnostics in Berlin)"
$ws.Range("D31").Value = "yes"
$ws.Range("E31").Value = "2025-06-19"

$ws.Range("A32").Value = 31
$ws.Range("B32").Value = "TARDBP"
$ws.Range("C32").Value = "Amyotrophic lateral sclerosis 10, with or without FTD (OMIM #612069; frequent diagnostics in Berlin)"
$ws.Range("D32").Value = "yes"
$ws.Range("E32").Value = "2025-06-19"

$ws.Range("A33").Value = 32
$ws.Range("B33").Value = "SOD1"
$ws.Range("C33").Value = "Amyotrophic lateral sclerosis 1 (OMIM #105400; frequent diagnostics in Berlin)"
$ws.Range("D33").Value = "yes"
$ws.Range("E33").Value = "2025-06-19"

$ws.Range("A34").Value = 33
$ws.Range("B34").Value = "MUC1"
$ws.Range("C34").Value = "Tubulointerstitial kidney disease, autosomal dominant, 2 (OMIM #174000; frequent diagnostics in Berlin)"
$ws.Range("D34").Value = "yes"
$ws.Range("E34").Value = "2025-06-19"

$ws.Range("A35").Value = 34
$ws.Range("B35").Value = "GBA"
$ws.Range("C35").Value = "Gaucher disease, type I (OMIM #230800; frequent diagnostics in Desden)"
$ws.Range("D35").Value = "yes"
$ws.Range("E35").Value = "2025-06-19"

$ws.Range("A36").Value = 35
$ws.Range("B36").Value = "LRRK2"
$ws.Range("C36").Value = "{Parkinson disease 8} (OMIM #607060; frequent diagnostics in Desden)"
$ws.Range("D36").Value = "yes"
$ws.Range("E36").Value = "2025-06-19"

$ws.Range("A37").Value = 36
$ws.Range("B37").Value = "GRN"
$ws.Range("C37").Value = "Frontotemporal dementia 2 (OMIM #607485; frequent diagnostics in Desden)"
$ws.Range("D37").Value = "yes"
$ws.Range("E37").Value = "2025-06-19"

$ws.Range("A38").Value = 37
$ws.Range("B38").Value = "MAPT"
$ws.Range("C38").Value = "Frontotemporal dementia 1, with or without parkinsonism (OMIM #600274; frequent diagnostics in Desden)"
$ws.Range("D38").Value = "yes"
$ws.Range("E38").Value = "2025-06-19"

$ws.Range("A39").Value = 38
$ws.Range("B39").Value = "PSEN1"
$ws.Range("C39").Value = "Dementia, frontotemporal (OMIM #600274; frequent diagnostics in Desden)"
$ws.Range("D39").Value = "yes"
$ws.Range("E39").Value = "2025-06-19"

$ws.Range("A40").Value = 39
$ws.Range("B40").Value = "PSEN2"
$ws.Range("C40").Value = "Alzheimer disease-4 (OMIM #606889; frequent diagnostics in Desden)"
$ws.Range("D40").Value = "yes"
$ws.Range("E40").Value = "2025-06-19"

$ws.Range("A41").Value = 40
$ws.Range("B41").Value = "APP"
$ws.Range("C41").Value = "Cerebral amyloid angiopathy, Dutch, Italian, Iowa, Flemish, Arctic variants (OMIM #605714; frequent diagnostics in Desden)"
$ws.Range("D41").Value = "yes"
$ws.Range("E41").Value = "2025-06-19"

$ws.Range("A42").Value = 41
$ws.Range("B42").Value = "DCTN1"
$ws.Range("C42").Value = "Neuronopathy, distal hereditary motor, autosomal dominant 14 (OMIM #607641; frequent diagnostics in Desden)"
$ws.Range("D42").Value = "yes"
$ws.Range("E42").Value = "2025-06-19"

# Re-fit column C (reason) now that it holds longer descriptions
$ws.Columns.Item(3).ColumnWidth = 114.3

# Highlight duplicate gene symbols (column B) the way Excel's
# "Highlight Cells Rules > Duplicate Values" gallery entry does
$dupRange = $ws.Range("B1:B1048576")
$dupFormat = $dupRange.FormatConditions.AddUniqueValues()
$dupFormat.DupeUnique = 1
$dupFormat.Font.Color = 393372
$dupFormat.Interior.Color = 13551615

